$wb = $excel.ActiveWorkbook

# --- Update existing sheets with recomputed values ---
$ws = $wb.Worksheets.Item("descriptives")
$ws.Range("C2").Value = [double]"787"
$ws.Range("D2").Value = [double]"0.006950631904356579"
$ws.Range("E2").Value = [double]"0.05229319073831001"
$ws.Range("F2").Value = [double]"97.85472485650796"
$ws.Range("G2").Value = [double]"86.37416634018692"
$ws.Range("H2").Value = [double]"11.48055851632104"
$ws.Range("C3").Value = [double]"421"
$ws.Range("D3").Value = [double]"0.002696208281190124"
$ws.Range("E3").Value = [double]"0.05376558133719576"
$ws.Range("F3").Value = [double]"95.10822366055712"
$ws.Range("G3").Value = [double]"90.56654012597475"
$ws.Range("H3").Value = [double]"4.541683534582362"

$ws = $wb.Worksheets.Item("coefficients")
$ws.Range("D2").Value = [double]"0.09480380187858262"
$ws.Range("E2").Value = [double]"0.01981462143092193"
$ws.Range("F2").Value = [double]"4.798949522832739"
$ws.Range("G2").Value = [double]"2.20074119876882e-05"
$ws.Range("H2").Value = [double]"0.05499871317828864"
$ws.Range("I2").Value = [double]"0.1343080343854703"
$ws.Range("J2").Value = [double]"40.38026157118141"
$ws.Range("D3").Value = [double]"0.06675628831344069"
$ws.Range("E3").Value = [double]"0.03003831250297576"
$ws.Range("F3").Value = [double]"2.225681572144141"
$ws.Range("G3").Value = [double]"0.05208371956087731"
$ws.Range("H3").Value = [double]"-0.0007413541287402573"
$ws.Range("I3").Value = [double]"0.133648414926698"
$ws.Range("J3").Value = [double]"9.320302586188001"
$ws.Range("D4").Value = [double]"0.2014534170335661"
$ws.Range("E4").Value = [double]"0.0689296459707017"
$ws.Range("F4").Value = [double]"2.963122571946534"
$ws.Range("G4").Value = [double]"0.01571249113603998"
$ws.Range("H4").Value = [double]"0.04851152953152397"
$ws.Range("I4").Value = [double]"0.3451650073861444"
$ws.Range("J4").Value = [double]"9.088961875362513"
$ws.Range("D5").Value = [double]"0.09101557961143948"
$ws.Range("E5").Value = [double]"0.05587487798336414"
$ws.Range("F5").Value = [double]"1.633438118306467"
$ws.Range("G5").Value = [double]"0.3497244726583097"
$ws.Range("H5").Value = [double]"-0.5502149122376441"
$ws.Range("I5").Value = [double]"0.6647214995021303"
$ws.Range("D6").Value = [double]"0.06938290263373881"
$ws.Range("E6").Value = [double]"0.07201046353820569"
$ws.Range("F6").Value = [double]"0.965061993022556"
$ws.Range("G6").Value = [double]"0.5113176494580161"
$ws.Range("H6").Value = [double]"-0.6887034325226536"
$ws.Range("I6").Value = [double]"0.7549962556010068"
$ws.Range("J6").Value = [double]"0.999999999999999"
$ws.Range("D7").Value = [double]"0.05617493319933421"
$ws.Range("E7").Value = [double]"0.01923951046584331"
$ws.Range("F7").Value = [double]"2.922846420761868"
$ws.Range("G7").Value = [double]"0.006735869609363888"
$ws.Range("H7").Value = [double]"0.01684833748883502"
$ws.Range("I7").Value = [double]"0.09532798985678863"
$ws.Range("J7").Value = [double]"28.41879720705413"
$ws.Range("D8").Value = [double]"0.003328434111490566"
$ws.Range("E8").Value = [double]"0.05613765566378266"
$ws.Range("F8").Value = [double]"0.05929079801319181"
$ws.Range("G8").Value = [double]"0.9599395979900474"
$ws.Range("H8").Value = [double]"-0.3457765397700591"
$ws.Range("I8").Value = [double]"0.3516239790544953"
$ws.Range("J8").Value = [double]"1.424607765223479"
$ws.Range("D9").Value = [double]"0.1016071824419475"
$ws.Range("E9").Value = [double]"0.06803836844033433"
$ws.Range("F9").Value = [double]"1.498551944941392"
$ws.Range("G9").Value = [double]"0.1823349398077589"
$ws.Range("H9").Value = [double]"-0.06254580648562731"
$ws.Range("I9").Value = [double]"0.2604076078408616"
$ws.Range("J9").Value = [double]"6.29905317958825"
$ws.Range("D10").Value = [double]"-0.03716279722333067"
$ws.Range("E10").Value = [double]"2.136512122426742e-17"
$ws.Range("F10").Value = [double]"-1740215709914876"
$ws.Range("G10").Value = [double]"3.658280802434097e-16"
$ws.Range("H10").Value = [double]"-0.03716279722333094"
$ws.Range("I10").Value = [double]"-0.0371627972233304"
$ws.Range("D11").Value = [double]"-0.0867829520916848"
$ws.Range("E11").Value = [double]"0.01778600421437634"
$ws.Range("F11").Value = [double]"-4.891587964268417"
$ws.Range("G11").Value = [double]"0.1283769302619888"
$ws.Range("H11").Value = [double]"-0.3031587723925541"
$ws.Range("I11").Value = [double]"0.1381026406794477"
$ws.Range("J11").Value = [double]"1"

$ws = $wb.Worksheets.Item("pairwise")
$ws.Range("C2").Value = [double]"0.7845951929272221"
$ws.Range("E2").Value = [double]"15.37448879295304"
$ws.Range("F2").Value = [double]"0.4446127837860866"
$ws.Range("G2").Value = [double]"0.8892255675721732"
$ws.Range("C3").Value = [double]"1.52197374291024"
$ws.Range("E3").Value = [double]"12.50958271745973"
$ws.Range("F3").Value = [double]"0.1528851567857516"
$ws.Range("G3").Value = [double]"0.7644257839287582"
$ws.Range("C4").Value = [double]"0.06445578880778553"
$ws.Range("E4").Value = [double]"1.103959442607108"
$ws.Range("F4").Value = [double]"0.9582560433218584"
$ws.Range("G4").Value = [double]"0.9776543772172571"
$ws.Range("C5").Value = [double]"0.3426949377405902"
$ws.Range("E5").Value = [double]"1.066789937473456"
$ws.Range("F5").Value = [double]"0.7870051312758293"
$ws.Range("G5").Value = [double]"0.9776543772172571"
$ws.Range("C6").Value = [double]"1.827244841532476"
$ws.Range("E6").Value = [double]"17.19680970945057"
$ws.Range("F6").Value = [double]"0.08506983293021878"
$ws.Range("G6").Value = [double]"0.7644257839287582"
$ws.Range("C7").Value = [double]"0.3848272713077692"
$ws.Range("E7").Value = [double]"1.365976419755909"
$ws.Range("F7").Value = [double]"0.7517174391515338"
$ws.Range("G7").Value = [double]"0.9776543772172571"
$ws.Range("C8").Value = [double]"0.03382073872765245"
$ws.Range("E8").Value = [double]"1.231644076509574"
$ws.Range("F8").Value = [double]"0.9776543772172571"
$ws.Range("G8").Value = [double]"0.9776543772172571"
$ws.Range("C9").Value = [double]"1.273264827496986"
$ws.Range("E9").Value = [double]"1.639235642331296"
$ws.Range("F9").Value = [double]"0.3539742264235008"
$ws.Range("G9").Value = [double]"0.8853777362822417"
$ws.Range("C10").Value = [double]"1.351802550463184"
$ws.Range("E10").Value = [double]"1.399792314785578"
$ws.Range("F10").Value = [double]"0.3541510945128967"
$ws.Range("G10").Value = [double]"0.8853777362822417"
$ws.Range("C11").Value = [double]"0.2388882099633046"
$ws.Range("E11").Value = [double]"1.912705596849554"
$ws.Range("F11").Value = [double]"0.8343335863579571"
$ws.Range("G11").Value = [double]"0.9776543772172571"
$ws.Range("C12").Value = [double]"0.8915232062478949"
$ws.Range("E12").Value = [double]"1.784209748033942"
$ws.Range("F12").Value = [double]"0.4763230766920143"
$ws.Range("G12").Value = [double]"0.570332065120611"
$ws.Range("C13").Value = [double]"0.6466878702976184"
$ws.Range("E13").Value = [double]"8.945989787947109"
$ws.Range("F13").Value = [double]"0.5340719865742932"
$ws.Range("G13").Value = [double]"0.570332065120611"
$ws.Range("C14").Value = [double]"4.855323843444371"
$ws.Range("E14").Value = [double]"28.41879720705413"
$ws.Range("F14").Value = [double]"3.976870764333308e-05"
$ws.Range("G14").Value = [double]"0.0003976870764333308"
$ws.Range("C15").Value = [double]"5.466778815403683"
$ws.Range("E15").Value = [double]"1.085278394104082"
$ws.Range("F15").Value = [double]"0.1013292331564061"
$ws.Range("G15").Value = [double]"0.3377641105213536"
$ws.Range("C16").Value = [double]"1.118158819698866"
$ws.Range("E16").Value = [double]"3.396672876258658"
$ws.Range("F16").Value = [double]"0.3362862038668719"
$ws.Range("G16").Value = [double]"0.4804088626669598"
$ws.Range("C17").Value = [double]"0.7215899118447767"
$ws.Range("E17").Value = [double]"1.424607765223478"
$ws.Range("F17").Value = [double]"0.570332065120611"
$ws.Range("G17").Value = [double]"0.570332065120611"
$ws.Range("C18").Value = [double]"1.533937097233646"
$ws.Range("E18").Value = [double]"1.677731893202403"
$ws.Range("F18").Value = [double]"0.2871574790881026"
$ws.Range("G18").Value = [double]"0.4785957984801711"
$ws.Range("C19").Value = [double]"2.045007135607972"
$ws.Range("E19").Value = [double]"6.299053179588256"
$ws.Range("F19").Value = [double]"0.08460543806470976"
$ws.Range("G19").Value = [double]"0.3377641105213536"
$ws.Range("C20").Value = [double]"2.686977262348441"
$ws.Range("E20").Value = [double]"1.459364654619941"
$ws.Range("F20").Value = [double]"0.1588864775911141"
$ws.Range("G20").Value = [double]"0.3972161939777852"
$ws.Range("C21").Value = [double]"2.80118479385572"
$ws.Range("E21").Value = [double]"0.9999999999999998"
$ws.Range("F21").Value = [double]"0.218290530984795"
$ws.Range("G21").Value = [double]"0.43658106196959"

# --- Add new "nr_studies" worksheet at the end ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ns = $wb.Worksheets.Add($null, $lastSheet)
$ns.Name = "nr_studies"

# Header row text
$ns.Range("A1").Value = "outcome"
$ns.Range("B1").Value = "country_id_europe_1_usa_north_america_2_asia_3_australia_4_south_america_5"
$ns.Range("C1").Value = "n_effect_sizes"
$ns.Range("D1").Value = "k_studies"
# Reuse the bold + centered header style already used on the other sheets
$headerSrc = $wb.Worksheets.Item("pairwise")
$headerSrc.Range("A1:D1").Copy() | Out-Null
$ns.Range("A1:D1").PasteSpecial(-4122) | Out-Null

$ns.Range("A2").Value = "NS"
$ns.Range("B2").NumberFormat = "@"
$ns.Range("B2").Value = "1"
$ns.Range("C2").Value = 517
$ns.Range("D2").Value = 48

$ns.Range("A3").Value = "NS"
$ns.Range("B3").NumberFormat = "@"
$ns.Range("B3").Value = "3"
$ns.Range("C3").Value = 59
$ns.Range("D3").Value = 11

$ns.Range("A4").Value = "NS"
$ns.Range("C4").Value = 58
$ns.Range("D4").Value = 6

$ns.Range("A5").Value = "NS"
$ns.Range("B5").NumberFormat = "@"
$ns.Range("B5").Value = "2"
$ns.Range("C5").Value = 177
$ns.Range("D5").Value = 13

$ns.Range("A6").Value = "NS"
$ns.Range("B6").NumberFormat = "@"
$ns.Range("B6").Value = "4"
$ns.Range("C6").Value = 27
$ns.Range("D6").Value = 2

$ns.Range("A7").Value = "NS"
$ns.Range("B7").NumberFormat = "@"
$ns.Range("B7").Value = "5"
$ns.Range("C7").Value = 7
$ns.Range("D7").Value = 2

$ns.Range("A8").Value = "NT"
$ns.Range("B8").NumberFormat = "@"
$ns.Range("B8").Value = "1"
$ns.Range("C8").Value = 282
$ns.Range("D8").Value = 35

$ns.Range("A9").Value = "NT"
$ns.Range("B9").NumberFormat = "@"
$ns.Range("B9").Value = "3"
$ns.Range("C9").Value = 63
$ns.Range("D9").Value = 8

$ns.Range("A10").Value = "NT"
$ns.Range("B10").NumberFormat = "@"
$ns.Range("B10").Value = "2"
$ns.Range("C10").Value = 51
$ns.Range("D10").Value = 3

$ns.Range("A11").Value = "NT"
$ns.Range("C11").Value = 17
$ns.Range("D11").Value = 2

$ns.Range("A12").Value = "NT"
$ns.Range("B12").NumberFormat = "@"
$ns.Range("B12").Value = "4"
$ns.Range("C12").Value = 18
$ns.Range("D12").Value = 1

$ns.Range("A13").Value = "NT"
$ns.Range("B13").NumberFormat = "@"
$ns.Range("B13").Value = "5"
$ns.Range("C13").Value = 7
$ns.Range("D13").Value = 2
